$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.026671460235407
$ws.Range("D2").Value = 1.038915501704312
$ws.Range("E2").Value = 1.026864470649486
$ws.Range("F2").Value = 1.038633994696341
$ws.Range("I2").Value = 1.038112585822517
$ws.Range("J2").Value = 1.031833601415962
$ws.Range("K2").Value = 1.041702175808341
$ws.Range("L2").Value = 1.029685872565378
$ws.Range("M2").Value = 1.0414214700446
$ws.Range("N2").Value = 1.033298923482952
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.027618517893215
$ws.Range("D3").Value = 1.039517420154822
$ws.Range("E3").Value = 1.027668378249994
$ws.Range("F3").Value = 1.039825051248679
$ws.Range("I3").Value = 1.038363400505145
$ws.Range("J3").Value = 1.032420530751993
$ws.Range("K3").Value = 1.042114755303684
$ws.Range("L3").Value = 1.030297394411474
$ws.Range("M3").Value = 1.042421574270376
$ws.Range("N3").Value = 1.033886686325951
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.02823153521702
$ws.Range("D4").Value = 1.039906678261499
$ws.Range("E4").Value = 1.028189122706892
$ws.Range("F4").Value = 1.040596144534266
$ws.Range("I4").Value = 1.038524289408417
$ws.Range("J4").Value = 1.032799923913078
$ws.Range("K4").Value = 1.042380803351435
$ws.Range("L4").Value = 1.030693002132187
$ws.Range("M4").Value = 1.043068541351774
$ws.Range("N4").Value = 1.034266618268842
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.028489296917604
$ws.Range("D5").Value = 1.040070267485894
$ws.Range("E5").Value = 1.028408177361583
$ws.Range("F5").Value = 1.040920408180215
$ws.Range("I5").Value = 1.038591590237115
$ws.Range("J5").Value = 1.032959327007715
$ws.Range("K5").Value = 1.042492429113739
$ws.Range("L5").Value = 1.030859293912296
$ws.Range("M5").Value = 1.04334048635163
$ws.Range("N5").Value = 1.034426247734156
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.028532579112883
$ws.Range("D6").Value = 1.040097731577098
$ws.Range("E6").Value = 1.028444965412981
$ws.Range("F6").Value = 1.040974859111872
$ws.Range("I6").Value = 1.03860287057343
$ws.Range("J6").Value = 1.032986085972073
$ws.Range("K6").Value = 1.042511158592662
$ws.Range("L6").Value = 1.030887213741135
$ws.Range("M6").Value = 1.043386144772575
$ws.Range("N6").Value = 1.034453044699314
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.028234979248713
$ws.Range("D7").Value = 1.039908864366155
$ws.Range("E7").Value = 1.028192049201846
$ws.Range("F7").Value = 1.040600476983323
$ws.Range("I7").Value = 1.038525190009684
$ws.Range("J7").Value = 1.032802054234285
$ws.Range("K7").Value = 1.042382295769275
$ws.Range("L7").Value = 1.030695224217242
$ws.Range("M7").Value = 1.0430721752511
$ws.Range("N7").Value = 1.034268751615349
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.026991479497743
$ws.Range("D8").Value = 1.039118968271485
$ws.Range("E8").Value = 1.027136037935612
$ws.Range("F8").Value = 1.039036435154124
$ws.Range("I8").Value = 1.038197640417926
$ws.Range("J8").Value = 1.032032037208373
$ws.Range("K8").Value = 1.041841798585023
$ws.Range("L8").Value = 1.029892557005587
$ws.Range("M8").Value = 1.041759494619606
$ws.Range("N8").Value = 1.033497641076938
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.024801886664735
$ws.Range("D9").Value = 1.037725421380781
$ws.Range("E9").Value = 1.025279561894678
$ws.Range("F9").Value = 1.036283445656545
$ws.Range("I9").Value = 1.037609718317114
$ws.Range("J9").Value = 1.030672211876908
$ws.Range("K9").Value = 1.040882380866303
$ws.Range("L9").Value = 1.028477512343915
$ws.Range("M9").Value = 1.039445099075232
$ws.Range("N9").Value = 1.032135884637529
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.023343269695345
$ws.Range("D10").Value = 1.036795363856408
$ws.Range("E10").Value = 1.02404489273704
$ws.Range("F10").Value = 1.034450151545215
$ws.Range("I10").Value = 1.037210580188302
$ws.Range("J10").Value = 1.029763706027299
$ws.Range("K10").Value = 1.040238123094334
$ws.Range("L10").Value = 1.02753375644222
$ws.Range("M10").Value = 1.037901296699984
$ws.Range("N10").Value = 1.031226088605406
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.022711941143666
$ws.Range("D11").Value = 1.036392412399755
$ws.Range("E11").Value = 1.023510986002651
$ws.Range("F11").Value = 1.033656791563418
$ws.Range("I11").Value = 1.037036049506107
$ws.Range("J11").Value = 1.029369855494115
$ws.Range("K11").Value = 1.039958061826366
$ws.Range("L11").Value = 1.027125014400329
$ws.Range("M11").Value = 1.037232603574449
$ws.Range("N11").Value = 1.030831678759289
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.022477477324276
$ws.Range("D12").Value = 1.036242705022422
$ws.Range("E12").Value = 1.023312777292145
$ws.Range("F12").Value = 1.033362172141674
$ws.Range("I12").Value = 1.0369709658749
$ws.Range("J12").Value = 1.02922349298759
$ws.Range("K12").Value = 1.039853871209441
$ws.Range("L12").Value = 1.026973176540826
$ws.Range("M12").Value = 1.036984188287201
$ws.Range("N12").Value = 1.030685108401217
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.022527768807139
$ws.Range("D13").Value = 1.036274819246729
$ws.Range("E13").Value = 1.023355288831039
$ws.Range("F13").Value = 1.03342536585848
$ws.Range("I13").Value = 1.036984938074629
$ws.Range("J13").Value = 1.029254891351876
$ws.Range("K13").Value = 1.03987622783023
$ws.Range("L13").Value = 1.02700574685
$ws.Range("M13").Value = 1.037037475703526
$ws.Range("N13").Value = 1.030716551354782
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.022692559479198
$ws.Range("D14").Value = 1.036380038212045
$ws.Range("E14").Value = 1.023494599800101
$ws.Range("F14").Value = 1.033632436809571
$ws.Range("I14").Value = 1.037030674875669
$ws.Range("J14").Value = 1.029357758526225
$ws.Range("K14").Value = 1.039949452725553
$ws.Range("L14").Value = 1.027112463694597
$ws.Range("M14").Value = 1.037212070162263
$ws.Range("N14").Value = 1.030819564612318
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.022794097723173
$ws.Range("D15").Value = 1.036444862721635
$ws.Range("E15").Value = 1.023580448226308
$ws.Range("F15").Value = 1.033760029288538
$ws.Range("I15").Value = 1.037058821029232
$ws.Range("J15").Value = 1.029421129273369
$ws.Range("K15").Value = 1.03999454737545
$ws.Range("L15").Value = 1.027178213777177
$ws.Range("M15").Value = 1.037319639202109
$ws.Range("N15").Value = 1.030883025353192
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.023385174443968
$ws.Range("D16").Value = 1.036822101678536
$ws.Range("E16").Value = 1.024080341469123
$ws.Range("F16").Value = 1.034502814054663
$ws.Range("I16").Value = 1.037222127408205
$ws.Range("J16").Value = 1.029789834905429
$ws.Range("K16").Value = 1.040256686876733
$ws.Range("L16").Value = 1.027560881469226
$ws.Range("M16").Value = 1.037945671068171
$ws.Range("N16").Value = 1.03125225458954
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.023756011683717
$ws.Range("D17").Value = 1.037058672978481
$ws.Range("E17").Value = 1.02439410284737
$ws.Range("F17").Value = 1.034968868324918
$ws.Range("I17").Value = 1.037324110017181
$ws.Range("J17").Value = 1.030020990969949
$ws.Range("K17").Value = 1.040420827947432
$ws.Range("L17").Value = 1.027800895249723
$ws.Range("M17").Value = 1.038338306098987
$ws.Range("N17").Value = 1.031483738922186
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.023972340050722
$ws.Range("D18").Value = 1.037196638693657
$ws.Range("E18").Value = 1.024577183301538
$ws.Range("F18").Value = 1.035240755134873
$ws.Range("I18").Value = 1.037383430574163
$ws.Range("J18").Value = 1.030155775832624
$ws.Range("K18").Value = 1.040516463100837
$ws.Range("L18").Value = 1.027940882551965
$ws.Range("M18").Value = 1.038567302603556
$ws.Range("N18").Value = 1.031618715194826
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.024046106711542
$ws.Range("D19").Value = 1.037243677609821
$ws.Range("E19").Value = 1.024639620639714
$ws.Range("F19").Value = 1.035333469201957
$ws.Range("I19").Value = 1.037403629491447
$ws.Range("J19").Value = 1.030201726415864
$ws.Range("K19").Value = 1.040549054285773
$ws.Range("L19").Value = 1.027988613115246
$ws.Range("M19").Value = 1.038645380965507
$ws.Range("N19").Value = 1.031664731033164
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.023716221740332
$ws.Range("D20").Value = 1.037033293415914
$ws.Range("E20").Value = 1.024360432103288
$ws.Range("F20").Value = 1.034918860437153
$ws.Range("I20").Value = 1.037313185226499
$ws.Range("J20").Value = 1.029996194717466
$ws.Range("K20").Value = 1.040403228075737
$ws.Range("L20").Value = 1.027775144935606
$ws.Range("M20").Value = 1.038296182236761
$ws.Range("N20").Value = 1.03145890745618
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.02264403163469
$ws.Range("D21").Value = 1.036349054756724
$ws.Range("E21").Value = 1.023453573202942
$ws.Range("F21").Value = 1.033571457653427
$ws.Range("I21").Value = 1.037017213568097
$ws.Range("J21").Value = 1.029327468605126
$ws.Range("K21").Value = 1.03992789432235
$ws.Range("L21").Value = 1.027081038594207
$ws.Range("M21").Value = 1.037160657356528
$ws.Range("N21").Value = 1.030789231676056
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.021970133095771
$ws.Range("D22").Value = 1.035918654958111
$ws.Range("E22").Value = 1.022884020439393
$ws.Range("F22").Value = 1.032724695599079
$ws.Range("I22").Value = 1.036829648344919
$ws.Range("J22").Value = 1.028906615648841
$ws.Range("K22").Value = 1.039628088220646
$ws.Range("L22").Value = 1.026644552237244
$ws.Range("M22").Value = 1.036446517345914
$ws.Range("N22").Value = 1.030367781060301
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.02232735755973
$ws.Range("D23").Value = 1.036146835765246
$ws.Range("E23").Value = 1.023185891560217
$ws.Range("F23").Value = 1.033173542126914
$ws.Range("I23").Value = 1.036929219992955
$ws.Range("J23").Value = 1.029129755367526
$ws.Range("K23").Value = 1.039787110381305
$ws.Range("L23").Value = 1.026875948748459
$ws.Range("M23").Value = 1.036825114649375
$ws.Range("N23").Value = 1.03059123766298
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.023734201024834
$ws.Range("D24").Value = 1.03704476141689
$ws.Range("E24").Value = 1.024375646250818
$ws.Range("F24").Value = 1.034941456709017
$ws.Range("I24").Value = 1.037318122176431
$ws.Range("J24").Value = 1.030007399214285
$ws.Range("K24").Value = 1.040411181026397
$ws.Range("L24").Value = 1.02778678042096
$ws.Range("M24").Value = 1.038315216261153
$ws.Range("N24").Value = 1.031470127864669
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.025367754978177
$ws.Range("D25").Value = 1.038085873553731
$ws.Range("E25").Value = 1.025758984050342
$ws.Range("F25").Value = 1.036994799792349
$ws.Range("I25").Value = 1.037762979723145
$ws.Range("J25").Value = 1.031024106384058
$ws.Range("K25").Value = 1.041131236272763
$ws.Range("L25").Value = 1.028843407401574
$ws.Range("M25").Value = 1.040043578888022
$ws.Range("N25").Value = 1.032488278875241
